# Fix a missing closing paren in the Consolas-styled code block on the
# "Example: Constraint Checking" slide:
#   if (!matchTypes(variable.type, expr)
# becomes
#   if (!matchTypes(variable.type, expr))

$p = $ppt.ActivePresentation

$needle = ", expr)"

# Locate the shape that contains the broken line instead of hard-coding
# slide/shape indices, so the script is resilient to reordering.
$targetShape = $null
$targetSlide = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                if ($tf.TextRange.Text.Contains($needle)) {
                    $targetShape = $shape
                    $targetSlide = $slide
                }
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the exact run of text ", expr)" that is missing the trailing ")".
$found = $tr.Find($needle, 0)

# Only touch the "expr)" portion (skip the leading ", ") so the existing
# ", " text/run is left completely untouched, matching the original
# formatting, and split a fresh run off for "expr))" .
$splitStart = $found.Start + 2
$splitLen = $found.Length - 2

$target = $tr.Characters($splitStart, $splitLen)
$target.Text = "expr))"
